function Set-TextValue($ws, $ref, $val) {
    $style = $ws.Range("A1").Style
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).Style = $style
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "246.02"
Set-TextValue $ws "E2" "-0.45%"
Set-TextValue $ws "G2" "9"

Set-TextValue $ws "D3" "29.99"
Set-TextValue $ws "E3" "-0.65%"
Set-TextValue $ws "G3" "9"

Set-TextValue $ws "D4" "5.151"
Set-TextValue $ws "E4" "-0.57%"
Set-TextValue $ws "G4" "9"

Set-TextValue $ws "E5" "0.34%"
Set-TextValue $ws "G5" "9"

Set-TextValue $ws "D6" "6.654"
Set-TextValue $ws "E6" "1.01%"
Set-TextValue $ws "G6" "9"

Set-TextValue $ws "D7" "3.237"
Set-TextValue $ws "E7" "5.18%"
Set-TextValue $ws "G7" "9"

Set-TextValue $ws "D8" "0.8519"
Set-TextValue $ws "E8" "-1.04%"
Set-TextValue $ws "G8" "9"

Set-TextValue $ws "D9" "0.8521"
Set-TextValue $ws "E9" "-3.22%"
Set-TextValue $ws "G9" "9"

Set-TextValue $ws "G10" "9"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws "D11" "0.07082"
Set-TextValue $ws "E11" "1.24%"
Set-TextValue $ws "G11" "9"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D12" "0.03259"
Set-TextValue $ws "E12" "11.59%"
Set-TextValue $ws "G12" "9"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws "D13" "0.09367"
Set-TextValue $ws "E13" "-0.19%"
Set-TextValue $ws "G13" "9"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D14" "0.001535"
Set-TextValue $ws "E14" "1.76%"
Set-TextValue $ws "G14" "9"

$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D15" "0.0005979"
Set-TextValue $ws "E15" "-0.54%"
Set-TextValue $ws "G15" "9"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D16" "0.005966"
Set-TextValue $ws "E16" "-2.71%"
Set-TextValue $ws "G16" "9"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D17" "3.511"
Set-TextValue $ws "E17" "0.20%"
Set-TextValue $ws "G17" "9"

$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws "D18" "2.204"
Set-TextValue $ws "E18" "-2.50%"
Set-TextValue $ws "G18" "9"

$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws "D19" "0.3166"
Set-TextValue $ws "E19" "0.68%"
Set-TextValue $ws "G19" "9"

$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws "D20" "0.03376"
Set-TextValue $ws "E20" "2.62%"
Set-TextValue $ws "G20" "9"

Set-TextValue $ws "E21" "-0.72%"
Set-TextValue $ws "G21" "9"

Set-TextValue $ws "D22" "3.488"
Set-TextValue $ws "E22" "-3.12%"
Set-TextValue $ws "G22" "9"

Set-TextValue $ws "E23" "2.42%"
Set-TextValue $ws "G23" "9"

Set-TextValue $ws "D24" "0.04124"
Set-TextValue $ws "E24" "0.02%"
Set-TextValue $ws "G24" "9"

Set-TextValue $ws "D25" "0.001228"
Set-TextValue $ws "E25" "1.37%"
Set-TextValue $ws "G25" "9"

Set-TextValue $ws "D26" "0.004142"
Set-TextValue $ws "E26" "-8.13%"
Set-TextValue $ws "G26" "9"

Set-TextValue $ws "E27" "1.78%"
Set-TextValue $ws "G27" "9"

Set-TextValue $ws "E28" "5.10%"
Set-TextValue $ws "G28" "9"

Set-TextValue $ws "G29" "9"

Set-TextValue $ws "G30" "9"

Set-TextValue $ws "G31" "9"

Set-TextValue $ws "G32" "9"

Set-TextValue $ws "G33" "9"

Set-TextValue $ws "G34" "9"

Set-TextValue $ws "G35" "9"

Set-TextValue $ws "G36" "9"

Set-TextValue $ws "G37" "9"

Set-TextValue $ws "G38" "9"

Set-TextValue $ws "G39" "9"

Set-TextValue $ws "E40" "-1.09%"
Set-TextValue $ws "G40" "9"

Set-TextValue $ws "D41" "0.005699"
Set-TextValue $ws "E41" "1.51%"
Set-TextValue $ws "G41" "9"

Set-TextValue $ws "D42" "0.1070"
Set-TextValue $ws "E42" "-0.14%"
Set-TextValue $ws "G42" "9"

Set-TextValue $ws "D43" "0.002200"
Set-TextValue $ws "E43" "0.10%"
Set-TextValue $ws "G43" "9"

Set-TextValue $ws "D44" "0.008474"
Set-TextValue $ws "E44" "-16.15%"
Set-TextValue $ws "G44" "9"

Set-TextValue $ws "D45" "0.00005482"
Set-TextValue $ws "E45" "7.26%"
Set-TextValue $ws "G45" "9"

Set-TextValue $ws "E46" "0.10%"
Set-TextValue $ws "G46" "9"

Set-TextValue $ws "D47" "0.07099"
Set-TextValue $ws "E47" "-20.14%"
Set-TextValue $ws "G47" "9"

Set-TextValue $ws "D48" "0.002229"
Set-TextValue $ws "E48" "-18.00%"
Set-TextValue $ws "G48" "9"

Set-TextValue $ws "D49" "0.00002100"
Set-TextValue $ws "E49" "0.10%"
Set-TextValue $ws "G49" "9"

Set-TextValue $ws "D50" "0.0002000"
Set-TextValue $ws "E50" "0.10%"
Set-TextValue $ws "G50" "9"

Set-TextValue $ws "G51" "9"
